$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds a text-formatted numeric string (inline string in the
# original workbook), e.g. "277.94" or "1.96%". Setting NumberFormat to "@" (Text)
# before assigning the value keeps Excel from reinterpreting these as numbers and
# losing formatting such as trailing zeros ("3.460") or the trailing "%" sign.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.10"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.54%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.911"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06414"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.951"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.65%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.246"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.74%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8827"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.52%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05014"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07522"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.58%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02881"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-8.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09008"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.31%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001582"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.43%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006423"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.53%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005862"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.99%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.316"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.31%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.00%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.64%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1337"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.05%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.916"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.30%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04427"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.54%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001175"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.35%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "13.85%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04140"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.82%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006820"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.03%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1175"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.06%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "13.79%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01173"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.87%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005205"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.13%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-22.20%"

Write-Output "Updated 58 price/volume cells on Sheet1"
